# Generate Report for Handback
# - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#   (this text is shared by the Overview sheet's zh-cn/de-de status columns and
#   the Status column on each language sheet).
# - Refreshed "Latest Handback DateTime" timestamps on the zh-cn / de-de sheets.
# - The stale "handback file is not the latest" warning is cleared now that the
#   handback is in sync.
# - A couple of columns that used to hold the (now gone) long warning text are
#   narrowed/widened to fit the new, shorter content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status cells for both rows ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de columns on the Overview sheet to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# Latest Handback DateTime refreshed.
$wsZhCn.Range("K2").Value = "2016-09-06 11:43:33"
$wsZhCn.Range("K3").Value = "2016-09-06 11:43:33"

# Error Detail no longer applies - the handback is now in sync.
$wsZhCn.Range("P2").Value = ""

# Column C (Status) needs to be wider for the new text; column P (Error Detail)
# can shrink now that it no longer holds the long warning message.
$wsZhCn.Columns.Item(3).ColumnWidth = 29.15
$wsZhCn.Columns.Item(16).ColumnWidth = 12.91

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Latest Handback DateTime refreshed.
$wsDeDe.Range("K2").Value = "2016-09-06 11:43:51"
$wsDeDe.Range("K3").Value = "2016-09-06 11:43:51"

# Error Detail no longer applies - the handback is now in sync.
$wsDeDe.Range("P2").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.15
$wsDeDe.Columns.Item(16).ColumnWidth = 12.91
